{"js": "const replacements = [\n  [\"2024-05-18 Saturday\", \"2024-05-19 Sunday\"],\n  [\"825\\u00F74=\", \"894\\u00F76=\"],\n  [\"469\\u00F77=\", \"390\\u00F74=\"],\n  [\"291\\u00F75=\", \"112\\u00F72=\"],\n  [\"656\\u00F78=\", \"738\\u00F72=\"],\n  [\"776\\u00F77=\", \"438\\u00F74=\"],\n  [\"252\\u00F79=\", \"739\\u00F75=\"],\n  [\"730\\u00F78=\", \"376\\u00F76=\"],\n  [\"824\\u00F74=\", \"720\\u00F77=\"],\n  [\"185\\u00F73=\", \"313\\u00F72=\"],\n  [\"630\\u00F77=\", \"974\\u00F74=\"],\n  [\"960\\u00F74=\", \"338\\u00F74=\"],\n  [\"696\\u00F76=\", \"637\\u00F73=\"],\n  [\"374\\u00F77=\", \"913\\u00F76=\"],\n  [\"357\\u00F76=\", \"315\\u00F77=\"],\n  [\"878\\u00F76=\", \"770\\u00F79=\"],\n  [\"551\\u00F76=\", \"683\\u00F73=\"],\n  [\"883\\u00F77=\", \"879\\u00F75=\"],\n  [\"573\\u00F72=\", \"526\\u00F79=\"],\n  [\"412\\u00F73=\", \"337\\u00F76=\"],\n  [\"819\\u00F73=\", \"400\\u00F79=\"],\n  [\"354\\u00F72=\", \"848\\u00F77=\"],\n  [\"442\\u00F72=\", \"801\\u00F74=\"],\n  [\"335\\u00F78=\", \"222\\u00F74=\"],\n  [\"524\\u00F73=\", \"141\\u00F75=\"],\n  [\"164\\u00F76=\", \"251\\u00F72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-05-18 Saturday\", \"2024-05-19 Sunday\"),\n    @(\"825\u00f74=\", \"894\u00f76=\"),\n    @(\"469\u00f77=\", \"390\u00f74=\"),\n    @(\"291\u00f75=\", \"112\u00f72=\"),\n    @(\"656\u00f78=\", \"738\u00f72=\"),\n    @(\"776\u00f77=\", \"438\u00f74=\"),\n    @(\"252\u00f79=\", \"739\u00f75=\"),\n    @(\"730\u00f78=\", \"376\u00f76=\"),\n    @(\"824\u00f74=\", \"720\u00f77=\"),\n    @(\"185\u00f73=\", \"313\u00f72=\"),\n    @(\"630\u00f77=\", \"974\u00f74=\"),\n    @(\"960\u00f74=\", \"338\u00f74=\"),\n    @(\"696\u00f76=\", \"637\u00f73=\"),\n    @(\"374\u00f77=\", \"913\u00f76=\"),\n    @(\"357\u00f76=\", \"315\u00f77=\"),\n    @(\"878\u00f76=\", \"770\u00f79=\"),\n    @(\"551\u00f76=\", \"683\u00f73=\"),\n    @(\"883\u00f77=\", \"879\u00f75=\"),\n    @(\"573\u00f72=\", \"526\u00f79=\"),\n    @(\"412\u00f73=\", \"337\u00f76=\"),\n    @(\"819\u00f73=\", \"400\u00f79=\"),\n    @(\"354\u00f72=\", \"848\u00f77=\"),\n    @(\"442\u00f72=\", \"801\u00f74=\"),\n    @(\"335\u00f78=\", \"222\u00f74=\"),\n    @(\"524\u00f73=\", \"141\u00f75=\"),\n    @(\"164\u00f76=\", \"251\u00f72=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $find\n    $rng.Find.Replacement.Text = $replace\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 1\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
